$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $val into the cell as plain text, even when it looks like a
# number (e.g. "1.001" or "27.015.39"), without leaving a lasting
# number-format change on the cell.
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# Row => Price (D), Volume/1h (E)
$updates = @{
    2  = @{ D = "27.015.39";      E = "  +0.39%  " }
    3  = @{ D = "1.887.99";       E = "  +1.40%  " }
    4  = @{ D = "1.001";          E = "  +0.09%  " }
    5  = @{ D = "306.08";         E = "  +0.42%  " }
    6  = @{ D = "0.9999";         E = "  +0.06%  " }
    7  = @{ D = "0.5176";         E = "  +2.78%  " }
    8  = @{ E = "  +2.95%  " }
    9  = @{ D = "0.07199";        E = "  +0.43%  " }
    10 = @{ D = "21.12";          E = "  +1.90%  " }
    11 = @{ D = "0.9009";         E = "  +0.67%  " }
    12 = @{ D = "0.07653";        E = "  +2.16%  " }
    13 = @{ D = "1.884.65";       E = "  +0.55%  " }
    14 = @{ D = "94.27";          E = "  -0.19%  " }
    15 = @{ E = "  +0.08%  " }
    16 = @{ E = "  +0.13%  " }
    17 = @{ D = "0.000008490";    E = "  -0.18%  " }
    18 = @{ D = "14.47";          E = "  +1.80%  " }
    19 = @{ D = "0.9999";         E = "  +0.10%  " }
    20 = @{ D = "27.062.13";      E = "  +0.37%  " }
    21 = @{ D = "5.047";          E = "  +0.41%  " }
    22 = @{ D = "2.117.91";       E = "  +0.51%  " }
    23 = @{ D = "10.59";          E = "  +2.08%  " }
    24 = @{ D = "6.383";          E = "  -0.48%  " }
    25 = @{ D = "2.307";          E = "  +10.58%  " }
    26 = @{ D = "146.22";         E = "  -0.94%  " }
    27 = @{ D = "18.03";          E = "  +0.91%  " }
    28 = @{ D = "1.726";          E = "  -2.94%  " }
    29 = @{ D = "114.24";         E = "  +1.05%  " }
    30 = @{ D = "4.916";          E = "  +5.15%  " }
    31 = @{ D = "4.788";          E = "  +1.94%  " }
    32 = @{ D = "0.09197";        E = "  -0.44%  " }
    33 = @{ D = "0.05037";        E = "  -2.02%  " }
    34 = @{ D = "1.239";          E = "  +7.32%  " }
    35 = @{ D = "0.7680";         E = "  +2.44%  " }
    36 = @{ D = "2.953";          E = "  -0.76%  " }
    37 = @{ D = "3.274";          E = "  +0.84%  " }
    38 = @{ D = "2.609";          E = "  +0.45%  " }
    39 = @{ D = "0.5591";         E = "  +0.51%  " }
    40 = @{ D = "0.01986";        E = "  -0.93%  " }
    41 = @{ D = "1.074";          E = "  +0.52%  " }
    42 = @{ D = "9.040";          E = "  +5.23%  " }
    43 = @{ D = "6.609";          E = "  +0.65%  " }
    44 = @{ D = "118.38";         E = "  +1.16%  " }
    45 = @{ D = "0.1505";         E = "  +2.19%  " }
    46 = @{ D = "0.4824";         E = "  +2.80%  " }
    47 = @{ D = "0.9995";         E = "  +0.10%  " }
    48 = @{ E = "  +0.90%  " }
    49 = @{ E = "  +2.40%  " }
    50 = @{ D = "37.66";          E = "  +2.53%  " }
    51 = @{ D = "63.96";          E = "  +1.39%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextValue $row 4 $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        Set-TextValue $row 5 $vals["E"]
    }
}
